$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Rule" label in B11 is changed from "R40" to the text "1".
# Assigning a numeric-looking string straight to .Value would make Excel
# store it as a real number, but the workbook needs it kept as text (a
# shared string), matching every other cell in that column. So the new
# text is staged in a scratch cell that is explicitly formatted as Text,
# and only its *value* is pasted into B11 - this swaps the cell's content
# without disturbing B11's existing style (s="23").
$scratch = $ws.Range("Z1")
$scratch.NumberFormat = "@"
$scratch.Value = "1"
$scratch.Copy()

$ws.Range("B11").PasteSpecial(-4163)  # xlPasteValues

# Clean up the scratch cell so it leaves no trace (contents + formatting).
$scratch.Clear()
$excel.CutCopyMode = $false
